# edit.ps1
# Updates country case/death statistics (data refresh) and reorders a few
# country rows per the commit "Update countries & provincias Spain".
#
# Three pairs/groups of countries effectively swap rank position because the
# underlying values changed: Belgica<->Suiza (rows 13/14), Brasil<->Portugal
# (rows 18/19) and a 5-way rotation among Macedonia/Bulgaria/Letonia/Libano/
# Eslovaquia (rows 79-83). The country name column (A) is updated alongside
# the numeric columns (B:H) for every row whose content changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp cell
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 13:22"

# Row 13: Suiza
$ws.Range("A13").Value = "Suiza"
$ws.Range("B13").Value = 22242
$ws.Range("C13").Value = 585
$ws.Range("D13").Value = 8056
$ws.Range("E13").Value = 13399
$ws.Range("F13").Value = 391
$ws.Range("G13").Value = 22
$ws.Range("H13").Value = 787

# Row 14: Belgica
$ws.Range("A14").Value = "Belgica"
$ws.Range("B14").Value = 22194
$ws.Range("C14").Value = 1380
$ws.Range("D14").Value = 4157
$ws.Range("E14").Value = 16002
$ws.Range("F14").Value = 1260
$ws.Range("G14").Value = 403
$ws.Range("H14").Value = 2035

# Row 17: Austria
$ws.Range("B17").Value = 12461
$ws.Range("C17").Value = 164
$ws.Range("E17").Value = 8172

# Row 18: Portugal
$ws.Range("A18").Value = "Portugal"
$ws.Range("B18").Value = 12442
$ws.Range("C18").Value = 712
$ws.Range("D18").Value = 184
$ws.Range("E18").Value = 11913
$ws.Range("F18").Value = 271
$ws.Range("G18").Value = 34
$ws.Range("H18").Value = 345

# Row 19: Brasil
$ws.Range("A19").Value = "Brasil"
$ws.Range("B19").Value = 12240
$ws.Range("C19").Value = 57
$ws.Range("D19").Value = 127
$ws.Range("E19").Value = 11547
$ws.Range("F19").Value = 296
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = 566

# Row 25: Noruega
$ws.Range("E25").Value = 5751
$ws.Range("G25").Value = 7
$ws.Range("H25").Value = 83

# Row 27: Dinamarca
$ws.Range("D27").Value = 1491
$ws.Range("E27").Value = 3284
$ws.Range("G27").Value = 16
$ws.Range("H27").Value = 203

# Row 28: India
$ws.Range("B28").Value = 4908
$ws.Range("C28").Value = 130
$ws.Range("E28").Value = 4389

# Row 32: Rumania
$ws.Range("E32").Value = 3773
$ws.Range("G32").Value = 8
$ws.Range("H32").Value = 184

# Row 33: Pakistan
$ws.Range("E33").Value = 3520
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = 55

# Row 79: Eslovaquia
$ws.Range("A79").Value = "Eslovaquia"
$ws.Range("B79").Value = 581
$ws.Range("C79").Value = 47
$ws.Range("D79").Value = 8
$ws.Range("E79").Value = 571
$ws.Range("F79").Value = 3
$ws.Range("H79").Value = 2

# Row 80: Republica de Macedonia
$ws.Range("A80").Value = "Republica de Macedonia"
$ws.Range("B80").Value = 570
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 30
$ws.Range("E80").Value = 517
$ws.Range("F80").Value = 15
$ws.Range("H80").Value = 23

# Row 81: Bulgaria
$ws.Range("A81").Value = "Bulgaria"
$ws.Range("B81").Value = 565
$ws.Range("C81").Value = 16
$ws.Range("D81").Value = 42
$ws.Range("E81").Value = 501
$ws.Range("F81").Value = 26
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 22

# Row 82: Letonia
$ws.Range("A82").Value = "Letonia"
$ws.Range("C82").Value = 6
$ws.Range("D82").Value = 16
$ws.Range("E82").Value = 530
$ws.Range("F82").Value = 5
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = 2

# Row 83: Libano
$ws.Range("A83").Value = "Libano"
$ws.Range("B83").Value = 548
$ws.Range("C83").Value = 7
$ws.Range("D83").Value = 60
$ws.Range("E83").Value = 469
$ws.Range("F83").Value = 27
$ws.Range("H83").Value = 19

# Row 85: Uzbekistan
$ws.Range("B85").Value = 504
$ws.Range("C85").Value = 47
$ws.Range("E85").Value = 472

# Row 113: Sri Lanka
$ws.Range("B113").Value = 183
$ws.Range("C113").Value = 5
$ws.Range("E113").Value = 135
